$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# Row 12: E12 gets a value and takes on the "yellow fill / thin border" look (like E9)
$ws.Range("E9").Copy()
$ws.Range("E12").PasteSpecial($xlPasteFormats)
$ws.Range("E12").Value = "1 revision: 1 insertion. 0 deletions"

# Row 13: D13 and E13 get values, both with the "yellow fill / thin border" look (like E9)
$ws.Range("E9").Copy()
$ws.Range("D13").PasteSpecial($xlPasteFormats)
$ws.Range("D13").Value = "1 revision: 1 insertion. 0 deletions"

$ws.Range("E9").Copy()
$ws.Range("E13").PasteSpecial($xlPasteFormats)
$ws.Range("E13").Value = "2 revisions: 2 insertions, 0 deletions"

# Row 14: D14 takes the "yellow fill / mixed border" look (like D9), E14 takes the "green fill" look (like E3)
$ws.Range("D9").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)
$ws.Range("D14").Value = "28 revisions: 23 insertions, 5 deletions"

$ws.Range("E3").Copy()
$ws.Range("E14").PasteSpecial($xlPasteFormats)
$ws.Range("E14").Value = "review complete - no change needed"

# Row 15: D15 takes the "yellow fill / mixed border" look (like D9), E15 takes the "yellow fill / thin border" look (like E9)
$ws.Range("D9").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Value = "39 revisions: 34 insertions, 5 deletions"

$ws.Range("E9").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("E15").Value = "8 revisions: 8 insertions, 0 deletions"

# Row 16: D16 takes the "yellow fill / mixed border" look (like D9); E16 is left as-is
$ws.Range("D9").Copy()
$ws.Range("D16").PasteSpecial($xlPasteFormats)
$ws.Range("D16").Value = "11 revisions: 9 insertions, 2 deletions"

$excel.CutCopyMode = $false

# Move the active selection to D16
$ws.Range("D16").Select()
